$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 500
$ws.Range("P3").Value = 300
$ws.Range("P4").Value = 200
$ws.Range("P5").Value = 0.6
$ws.Range("P6").Value = 180
$ws.Range("P7").Value = 90
$ws.Range("P8").Value = 90
$ws.Range("P9").Value = 0.5
$ws.Range("P11").Value = 320
$ws.Range("P12").Value = 290
$ws.Range("P13").Value = 0
$ws.Range("P14").Value = 30
$ws.Range("P15").Value = 300
$ws.Range("P16").Value = 60
$ws.Range("P17").Value = 50
$ws.Range("P18").Value = 40
$ws.Range("P19").Value = 150

$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("P21").Select()
